# Auto-generated edit script: updates column F (想去人数 / interest count)
# values across the four worksheets, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 2486
$ws.Cells.Item(7, 6).Value = 57
$ws.Cells.Item(8, 6).Value = 1835
$ws.Cells.Item(9, 6).Value = 3155
$ws.Cells.Item(10, 6).Value = 195
$ws.Cells.Item(11, 6).Value = 4649
$ws.Cells.Item(12, 6).Value = 429
$ws.Cells.Item(13, 6).Value = 250
$ws.Cells.Item(14, 6).Value = 144
$ws.Cells.Item(15, 6).Value = 594
$ws.Cells.Item(17, 6).Value = 4
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(19, 6).Value = 632
$ws.Cells.Item(20, 6).Value = 271
$ws.Cells.Item(21, 6).Value = 12
$ws.Cells.Item(24, 6).Value = 322
$ws.Cells.Item(25, 6).Value = 4632
$ws.Cells.Item(26, 6).Value = 9
$ws.Cells.Item(29, 6).Value = 5089
$ws.Cells.Item(30, 6).Value = 10
$ws.Cells.Item(31, 6).Value = 1164
$ws.Cells.Item(33, 6).Value = 634
$ws.Cells.Item(34, 6).Value = 4395
$ws.Cells.Item(36, 6).Value = 60
$ws.Cells.Item(37, 6).Value = 107
$ws.Cells.Item(38, 6).Value = 744
$ws.Cells.Item(39, 6).Value = 46
$ws.Cells.Item(40, 6).Value = 680
$ws.Cells.Item(41, 6).Value = 677

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 6
$ws.Cells.Item(6, 6).Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 1062

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1062
$ws.Cells.Item(9, 6).Value = 2486
$ws.Cells.Item(10, 6).Value = 57
$ws.Cells.Item(11, 6).Value = 1835
$ws.Cells.Item(13, 6).Value = 3155
$ws.Cells.Item(14, 6).Value = 195
$ws.Cells.Item(15, 6).Value = 4649
$ws.Cells.Item(16, 6).Value = 429
$ws.Cells.Item(17, 6).Value = 250
$ws.Cells.Item(18, 6).Value = 144
$ws.Cells.Item(19, 6).Value = 594
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(22, 6).Value = 4
$ws.Cells.Item(23, 6).Value = 632
$ws.Cells.Item(24, 6).Value = 271
$ws.Cells.Item(25, 6).Value = 12
$ws.Cells.Item(29, 6).Value = 322
$ws.Cells.Item(30, 6).Value = 4632
$ws.Cells.Item(31, 6).Value = 9
$ws.Cells.Item(34, 6).Value = 5089
$ws.Cells.Item(35, 6).Value = 10
$ws.Cells.Item(36, 6).Value = 1164
$ws.Cells.Item(38, 6).Value = 634
$ws.Cells.Item(39, 6).Value = 4395
$ws.Cells.Item(41, 6).Value = 6
$ws.Cells.Item(42, 6).Value = 60
$ws.Cells.Item(43, 6).Value = 107
$ws.Cells.Item(44, 6).Value = 744
$ws.Cells.Item(45, 6).Value = 46
$ws.Cells.Item(46, 6).Value = 680
$ws.Cells.Item(47, 6).Value = 677
$ws.Cells.Item(49, 6).Value = 40
